# Update gh-pages to output generated at 456a3b4
# Apply the same data refresh to both the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets, which contain duplicated rows.

$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Row 4: "最低票价" becomes unavailable -> text "不可售"
    $ws.Range("G4").Value = "不可售"

    # Row 5: "想去人数" count update
    $ws.Range("F5").Value = 5054

    # Row 6: "想去人数" count update
    $ws.Range("F6").Value = 375

    # Row 8: "想去人数" count update
    $ws.Range("F8").Value = 296

    # Row 9: "想去人数" count update
    $ws.Range("F9").Value = 762
}
